# edit.ps1 — apply the changes described by the diff:
#   1) Bump the cached "datetimeFigureOut" auto-date field text from
#      2017/3/28 -> 2017/4/2 everywhere it is cached (slide master,
#      every slide layout, and the notes master).
#   2) On slide 9, split the single run "是否登陸" into two runs
#      "是否" + "登入" (same rPr), i.e. change the visible text to
#      "是否登入".

$p = $ppt.ActivePresentation

$oldDate = "2017/3/28"
$newDate = "2017/4/2"

# --- helper: walk a Shapes collection and fix any date placeholder ---
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# 1) Slide master
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# 2) Every slide layout hanging off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# 3) Notes master
$notesMaster = $p.NotesMaster
Update-DatePlaceholders $notesMaster.Shapes

# --- slide 9: split "是否登陸" into "是否" + "登入" ---
$slide9 = $p.Slides.Item(9)
$shape = $slide9.Shapes.Item("文字方塊 8")
$tr = $shape.TextFrame.TextRange
if ($tr.Text -eq "是否登陸") {
    $tr.Text = "是否"
    $tr.InsertAfter("登入")
}
